$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New goal rows to add below "introduction" / "goodbye".
# Shared-string table indices are allocated in the order values are first
# written, so the JSON-blob text is written before its short label to match
# the target ordering (…JSON, label, label, label, label, JSON).
$ws.Range("D4").Value = "      ""basicVerification"": {
        ""achievementConditions"": ""This goal covers the basic information that the Customer Service Representative must verify in every call. A piece of information is verified if the borrower mentions it at any point, even if the Customer Service Representative does not repeat it. This goal is achieved if the Customer Service Representative verifies all 'borrowerInformation.'"",
        ""borrowerInformation"": [
          ""The borrower's first name"",
          ""The borrower's last name"",
          ""The borrower's state""
        ]
      }"
$ws.Range("A4").Value = "basicVerification"

$ws.Range("A5").Value = "fullVerification"

$ws.Range("A6").Value = "myVUVerification"

$ws.Range("A7").Value = "bankingVerification"
$ws.Range("D7").Value = "      ""bankingVerification"": {
        ""achievementConditions"": ""The Customer Service Representative verifies all of the borrower's 'requiredInformation' and at least one of the borrower's 'otherInformation.' A piece of information is verified if the borrower mentions it at any point."",
        ""requiredInformation"": [
          ""Full property address, including city, state, and ZIP code"",
          ""At least the last four digits of the borrower's Social Security Number"",
          ""The customer's first and last name"",
          ""The customer's Neighbors Bank account number""
        ],
        ""otherInformation"": [
          ""The joint owner on one of the customer's accounts"",
          ""The amount of a recent transaction using the customer's accounts"",
          ""The customer's email address on file"",
          ""The beneficiary of the customer's account"",
          ""The customer's current account balance"",
          ""The date of the last transaction using the customer's accoount""
        ]
      }"

# Match row height + wrap-text style used by the "introduction" row (row 2)
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(7).RowHeight = 409.5
$ws.Range("D4").WrapText = $true
$ws.Range("D7").WrapText = $true

# Update view state: scroll so row 2 is at top, select D5
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D5").Select()
